# Production doc builds - 2021/02/18 17:54:26 UTC
#
# Resizes the "AWS Cloud" bounding rectangle, groups three icon+label pairs
# (AWS IAM, AWS Systems Manager, Amazon EventBridge) into their own group
# shapes, removes the "Amazon EC2" icon+label, and repositions the
# "AWS Secrets Manager" icon+label.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU per point (914400 EMU per inch / 72 points per inch).
$EMU_PER_PT = 12700.0

function EmuToPt($emu) {
    # +0.5 nudges the float conversion so it round-trips to the exact EMU
    # value once PowerPoint re-quantizes points back to EMU internally.
    return ([double]$emu + 0.5) / $EMU_PER_PT
}

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

function Get-ShapeIndexById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $i
        }
    }
    return -1
}

# 1. Shrink the "AWS Cloud" outline rectangle.
$rect = Get-ShapeById $s.Shapes 25
$rect.Width = EmuToPt 8034612
$rect.Height = EmuToPt 3925693

# 2. Group the "AWS IAM" icon (30) + label (31) and move the group.
$idxIamIcon = Get-ShapeIndexById $s.Shapes 30
$idxIamLabel = Get-ShapeIndexById $s.Shapes 31
$iamGroup = $s.Shapes.Range(@($idxIamIcon, $idxIamLabel)).Group()
$iamGroup.Name = "Group 9"
$iamGroup.Left = EmuToPt 1297526
$iamGroup.Top = EmuToPt 2910805

# 3. Remove the "Amazon EC2" icon (36) + label (37) entirely.
$ec2Icon = Get-ShapeById $s.Shapes 36
$ec2Icon.Delete()
$ec2Label = Get-ShapeById $s.Shapes 37
$ec2Label.Delete()

# 4. Group the "AWS Systems Manager" icon (34) + label (35) and move it.
$idxSsmIcon = Get-ShapeIndexById $s.Shapes 34
$idxSsmLabel = Get-ShapeIndexById $s.Shapes 35
$ssmGroup = $s.Shapes.Range(@($idxSsmIcon, $idxSsmLabel)).Group()
$ssmGroup.Name = "Group 11"
$ssmGroup.Left = EmuToPt 4759829
$ssmGroup.Top = EmuToPt 2909075

# 5. Group the "Amazon EventBridge" icon (40) + label (41) and move it.
$idxEvbIcon = Get-ShapeIndexById $s.Shapes 40
$idxEvbLabel = Get-ShapeIndexById $s.Shapes 41
$evbGroup = $s.Shapes.Range(@($idxEvbIcon, $idxEvbLabel)).Group()
$evbGroup.Name = "Group 10"
$evbGroup.Left = EmuToPt 2834433
$evbGroup.Top = EmuToPt 2910805

# 6. Reposition the "AWS Secrets Manager" icon (44) and label (45); these
#    stay ungrouped, just move to their new spot.
$secretsIcon = Get-ShapeById $s.Shapes 44
$secretsIcon.Left = EmuToPt 7234289
$secretsIcon.Top = EmuToPt 2914580

$secretsLabel = Get-ShapeById $s.Shapes 45
$secretsLabel.Left = EmuToPt 6469114
$secretsLabel.Top = EmuToPt 3658548
